$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-apply the text number format to the header/label cells (same visual
# format as before -- Text "@" -- mirrors the style refresh seen in the
# authoring tool's re-export).
$ws.Range("A1:C1").NumberFormat = "@"
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A3").NumberFormat = "@"

# Updated prediction values for the second column.
$ws.Range("B2").Value = 10990.46771063232
$ws.Range("B3").Value = 7895.4244621250637
